$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 1): swap the Chinese labels for the English ones used by
# the refactored demo, and give A1/B1 the same header style C1 already had.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "Date"

$ws.Range("C1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

# Placeholder row (row 2): reorder so the templates line up with their
# matching header column (name/number/date).
$ws.Range("A2").Value = "{.name}"
$ws.Range("B2").Value = "{.number}"
$ws.Range("C2").Value = "{.date}"

# Match the author's final selection.
$ws.Range("D8").Select()
